$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 39
$newRow  = 40

# Every row in this log sheet shares the same cell style (centered
# horizontal/vertical alignment). Copy the formatting from the previous
# (last existing) row onto the new row first, so the freshly appended
# row - including its blank cells (F/H) - keeps that same style once
# saved, instead of those blank cells being dropped entirely.
$ws.Range("A$lastRow`:H$lastRow").Copy()
$ws.Range("A$newRow`:H$newRow").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Cells.Item($newRow, 1).Value = "2025-08-21 09:39:57 UTC"
$ws.Cells.Item($newRow, 2).Value = "2025-08-21 15:09:57 IST"
$ws.Cells.Item($newRow, 3).Value = "SKIPPED"
$ws.Cells.Item($newRow, 4).Value = "No change in PDF. Skipping download & Excel update."
$ws.Cells.Item($newRow, 5).Value = "https://nalcoindia.com/wp-content/uploads/2019/01/INGOT-21-08-2025.pdf"
$ws.Cells.Item($newRow, 6).Value = ""
$ws.Cells.Item($newRow, 7).Value = 0
$ws.Cells.Item($newRow, 8).Value = ""
